$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "25.624.77"
Set-TextValue $ws "E2" "  -4.16%  "
Set-TextValue $ws "D3" "1.811.93"
Set-TextValue $ws "E3" "  -2.92%  "
Set-TextValue $ws "D4" "1.001"
Set-TextValue $ws "E4" "  -0.09%  "
Set-TextValue $ws "D5" "276.58"
Set-TextValue $ws "E5" "  -8.07%  "
Set-TextValue $ws "D6" "1.001"
Set-TextValue $ws "E6" "  -0.07%  "
Set-TextValue $ws "D7" "0.5035"
Set-TextValue $ws "E7" "  -5.44%  "
Set-TextValue $ws "D8" "0.3491"
Set-TextValue $ws "E8" "  -6.40%  "
Set-TextValue $ws "D9" "44.29"
Set-TextValue $ws "E9" "  -2.23%  "
Set-TextValue $ws "D10" "0.06654"
Set-TextValue $ws "E10" "  -7.17%  "
Set-TextValue $ws "D11" "19.58"
Set-TextValue $ws "E11" "  -9.14%  "
Set-TextValue $ws "D12" "0.8102"
Set-TextValue $ws "E12" "  -8.78%  "
Set-TextValue $ws "D13" "0.07882"
Set-TextValue $ws "E13" "  -3.65%  "
Set-TextValue $ws "D14" "1.831.26"
Set-TextValue $ws "E14" "  -1.99%  "
Set-TextValue $ws "D15" "5.049"
Set-TextValue $ws "E15" "  -4.69%  "
Set-TextValue $ws "D16" "86.89"
Set-TextValue $ws "E16" "  -6.15%  "
Set-TextValue $ws "E17" "  -0.11%  "
Set-TextValue $ws "D18" "13.98"
Set-TextValue $ws "D19" "1.003"
Set-TextValue $ws "E19" "  +0.14%  "
Set-TextValue $ws "D20" "0.000007951"
Set-TextValue $ws "E20" "  -6.18%  "
Set-TextValue $ws "D21" "25.668.72"
Set-TextValue $ws "E21" "  -4.21%  "
Set-TextValue $ws "D22" "4.729"
Set-TextValue $ws "E22" "  -4.98%  "
Set-TextValue $ws "D23" "9.919"
Set-TextValue $ws "E23" "  -6.48%  "
Set-TextValue $ws "D24" "6.100"
Set-TextValue $ws "E24" "  -4.00%  "
Set-TextValue $ws "D25" "2.245"
Set-TextValue $ws "E25" "  -2.05%  "
Set-TextValue $ws "D26" "142.52"
Set-TextValue $ws "E26" "  -2.15%  "
Set-TextValue $ws "D27" "1.659"
Set-TextValue $ws "E27" "  -3.71%  "
Set-TextValue $ws "D28" "17.07"
Set-TextValue $ws "E28" "  -5.22%  "
Set-TextValue $ws "D29" "108.44"
Set-TextValue $ws "E29" "  -4.55%  "
Set-TextValue $ws "D30" "4.278"
Set-TextValue $ws "E30" "  -8.60%  "
Set-TextValue $ws "D31" "4.202"
Set-TextValue $ws "E31" "  -8.76%  "
Set-TextValue $ws "D32" "0.08743"
Set-TextValue $ws "E32" "  -4.12%  "
Set-TextValue $ws "D33" "0.04814"
Set-TextValue $ws "E33" "  -3.89%  "
Set-TextValue $ws "B34" "ImmutableX"
Set-TextValue $ws "C34" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D34" "0.7184"
Set-TextValue $ws "E34" "  -10.45%  "
Set-TextValue $ws "B35" "HuobiToken"
Set-TextValue $ws "C35" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D35" "2.871"
Set-TextValue $ws "E35" "  -2.61%  "
Set-TextValue $ws "D36" "1.123"
Set-TextValue $ws "E36" "  -4.18%  "
Set-TextValue $ws "B37" "Frax"
Set-TextValue $ws "C37" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws "D37" "1.002"
Set-TextValue $ws "E37" "  +0.09%  "
Set-TextValue $ws "B38" "MXToken"
Set-TextValue $ws "C38" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D38" "3.132"
Set-TextValue $ws "E38" "  -1.70%  "
Set-TextValue $ws "B39" "RenderToken"
Set-TextValue $ws "C39" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D39" "2.329"
Set-TextValue $ws "E39" "  -13.06%  "
Set-TextValue $ws "B40" "VeChain"
Set-TextValue $ws "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D40" "0.01839"
Set-TextValue $ws "E40" "  -5.52%  "
Set-TextValue $ws "B41" "TheSandbox"
Set-TextValue $ws "C41" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws "D41" "0.5067"
Set-TextValue $ws "E41" "  -16.44%  "
Set-TextValue $ws "B42" "TrustWalletToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D42" "0.9486"
Set-TextValue $ws "E42" "  -11.05%  "
Set-TextValue $ws "B43" "Quant"
Set-TextValue $ws "C43" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws "D43" "114.76"
Set-TextValue $ws "E43" "  +0.27%  "
Set-TextValue $ws "B44" "FraxShare"
Set-TextValue $ws "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D44" "6.186"
Set-TextValue $ws "E44" "  -4.85%  "
Set-TextValue $ws "B45" "Aptos"
Set-TextValue $ws "C45" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D45" "7.853"
Set-TextValue $ws "E45" "  -10.21%  "
Set-TextValue $ws "B46" "PaxDollar"
Set-TextValue $ws "C46" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws "D46" "1.002"
Set-TextValue $ws "E46" "  +0.10%  "
Set-TextValue $ws "D47" "0.1355"
Set-TextValue $ws "E47" "  -9.06%  "
Set-TextValue $ws "B48" "Decentraland"
Set-TextValue $ws "C48" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws "D48" "0.4459"
Set-TextValue $ws "E48" "  -14.44%  "
Set-TextValue $ws "B49" "Elrond"
Set-TextValue $ws "C49" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws "D49" "36.23"
Set-TextValue $ws "E49" "  -3.27%  "
Set-TextValue $ws "B50" "EnergySwap"
Set-TextValue $ws "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D50" "9.165"
Set-TextValue $ws "E50" "  -7.86%  "
Set-TextValue $ws "B51" "NEARProtocol"
Set-TextValue $ws "C51" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D51" "1.481"
Set-TextValue $ws "E51" "  -9.77%  "
